# Applies numeric updates to several cells across multiple worksheets
# as described by the source diff (profit/price recalculation updates).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3050
$ws.Range("I19").Value = 3500
$ws.Range("J19").Value = 2900
$ws.Range("K19").Value = 3500
$ws.Range("L19").Value = 2900
$ws.Range("M19").Value = -3325
$ws.Range("N19").Value = -3250
$ws.Range("H64").Value = 55168.74
$ws.Range("I64").Value = 85691.664
$ws.Range("J64").Value = 2843.7144
$ws.Range("K64").Value = 85691.664
$ws.Range("L64").Value = 2843.7144
$ws.Range("M64").Value = -85443.664
$ws.Range("N64").Value = -3339.7144
$ws.Range("H67").Value = 55168.74
$ws.Range("I67").Value = 85691.664
$ws.Range("J67").Value = 2843.7144
$ws.Range("K67").Value = 85691.664
$ws.Range("L67").Value = 2843.7144
$ws.Range("M67").Value = -84833.664
$ws.Range("N67").Value = -4559.7144
$ws.Range("H97").Value = 51224
$ws.Range("I97").Value = 800
$ws.Range("J97").Value = 56826.668
$ws.Range("K97").Value = 2400
$ws.Range("L97").Value = 170480.004
$ws.Range("M97").Value = -1904
$ws.Range("N97").Value = -171472.004
$ws.Range("H100").Value = 919
$ws.Range("I100").Value = 1000.3333
$ws.Range("J100").Value = 864.7778
$ws.Range("K100").Value = 1000.3333
$ws.Range("L100").Value = 864.7778
$ws.Range("M100").Value = -459.3333
$ws.Range("N100").Value = -1946.7778
$ws.Range("H112").Value = 1335.92
$ws.Range("J112").Value = 1489.3334
$ws.Range("L112").Value = 4468.0002
$ws.Range("N112").Value = -6684.0002
$ws.Range("H129").Value = 2428.3677
$ws.Range("J129").Value = 1210.0652
$ws.Range("L129").Value = 3630.1956
$ws.Range("N129").Value = -13630.1956
$ws.Range("H137").Value = 1059.3518
$ws.Range("I137").Value = 1044.58
$ws.Range("K137").Value = 3133.74
$ws.Range("M137").Value = -583.7399999999998
$ws.Range("H141").Value = 1852.4222
$ws.Range("I141").Value = 1689.738
$ws.Range("J141").Value = 4130
$ws.Range("K141").Value = 5069.214
$ws.Range("L141").Value = 12390
$ws.Range("M141").Value = 110.7860000000001
$ws.Range("N141").Value = -22750

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 34653.22
$ws.Range("I32").Value = 7435.6875
$ws.Range("K32").Value = 7435.6875
$ws.Range("M32").Value = -7148.6875
$ws.Range("H61").Value = 1279.875
$ws.Range("I61").Value = 1067.6571
$ws.Range("J61").Value = 2765.4
$ws.Range("K61").Value = 1067.6571
$ws.Range("L61").Value = 2765.4
$ws.Range("M61").Value = -855.6570999999999
$ws.Range("N61").Value = -3189.4
$ws.Range("H74").Value = 469.66666
$ws.Range("I74").Value = 434.5
$ws.Range("J74").Value = 891.6667
$ws.Range("K74").Value = 434.5
$ws.Range("L74").Value = 891.6667
$ws.Range("M74").Value = 439.5
$ws.Range("N74").Value = -2639.6667
$ws.Range("H77").Value = 469.66666
$ws.Range("I77").Value = 434.5
$ws.Range("J77").Value = 891.6667
$ws.Range("K77").Value = 2172.5
$ws.Range("L77").Value = 4458.3335
$ws.Range("M77").Value = 2195.5
$ws.Range("N77").Value = -13194.3335
$ws.Range("H110").Value = 83508936
$ws.Range("I110").Value = 100210536
$ws.Range("J110").Value = 950
$ws.Range("K110").Value = 100210536
$ws.Range("L110").Value = 950
$ws.Range("M110").Value = -100208491
$ws.Range("N110").Value = -5040
$ws.Range("H136").Value = 1279.875
$ws.Range("I136").Value = 1067.6571
$ws.Range("J136").Value = 2765.4
$ws.Range("K136").Value = 3202.9713
$ws.Range("L136").Value = 8296.200000000001
$ws.Range("M136").Value = -652.9712999999997
$ws.Range("N136").Value = -13396.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2071.8948
$ws.Range("I134").Value = 1870.8628
$ws.Range("J134").Value = 3780.6667
$ws.Range("K134").Value = 5612.588400000001
$ws.Range("L134").Value = 11342.0001
$ws.Range("M134").Value = -3077.588400000001
$ws.Range("N134").Value = -16412.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1426
$ws.Range("I16").Value = 966.4286
$ws.Range("K16").Value = 966.4286
$ws.Range("M16").Value = -679.4286
$ws.Range("H31").Value = 24972.016
$ws.Range("I31").Value = 1454.7028
$ws.Range("J31").Value = 59777.64
$ws.Range("K31").Value = 1454.7028
$ws.Range("L31").Value = 59777.64
$ws.Range("M31").Value = -1159.7028
$ws.Range("N31").Value = -60367.64
$ws.Range("H34").Value = 24972.016
$ws.Range("I34").Value = 1454.7028
$ws.Range("J34").Value = 59777.64
$ws.Range("K34").Value = 1454.7028
$ws.Range("L34").Value = 59777.64
$ws.Range("M34").Value = -1252.7028
$ws.Range("N34").Value = -60181.64
$ws.Range("H62").Value = 2417.8333
$ws.Range("I62").Value = 2162.8
$ws.Range("J62").Value = 2600
$ws.Range("K62").Value = 2162.8
$ws.Range("L62").Value = 2600
$ws.Range("M62").Value = -1538.8
$ws.Range("N62").Value = -3848
$ws.Range("H65").Value = 2417.8333
$ws.Range("I65").Value = 2162.8
$ws.Range("J65").Value = 2600
$ws.Range("K65").Value = 10814
$ws.Range("L65").Value = 13000
$ws.Range("M65").Value = -7694
$ws.Range("N65").Value = -19240
$ws.Range("H99").Value = 38998
$ws.Range("I99").Value = 8490
$ws.Range("J99").Value = 100014
$ws.Range("K99").Value = 8490
$ws.Range("L99").Value = 100014
$ws.Range("M99").Value = -6992
$ws.Range("N99").Value = -103010
$ws.Range("H107").Value = 8167.9287
$ws.Range("I107").Value = 15382.857
$ws.Range("J107").Value = 953
$ws.Range("K107").Value = 15382.857
$ws.Range("L107").Value = 953
$ws.Range("M107").Value = -13462.857
$ws.Range("N107").Value = -4793
$ws.Range("H113").Value = 1426
$ws.Range("I113").Value = 966.4286
$ws.Range("K113").Value = 966.4286
$ws.Range("M113").Value = 1203.5714
$ws.Range("H126").Value = 38998
$ws.Range("I126").Value = 8490
$ws.Range("J126").Value = 100014
$ws.Range("K126").Value = 25470
$ws.Range("L126").Value = 300042
$ws.Range("M126").Value = -23000
$ws.Range("N126").Value = -304982
$ws.Range("H132").Value = 4161.091
$ws.Range("I132").Value = 3792.72
$ws.Range("K132").Value = 11378.16
$ws.Range("M132").Value = -8848.16

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 915.55884
$ws.Range("I5").Value = 617.8421
$ws.Range("J5").Value = 1292.6666
$ws.Range("K5").Value = 1853.5263
$ws.Range("L5").Value = 3877.9998
$ws.Range("M5").Value = -1741.5263
$ws.Range("N5").Value = -4101.9998
$ws.Range("H131").Value = 7987.217
$ws.Range("J131").Value = 8023.6465
$ws.Range("L131").Value = 24070.9395
$ws.Range("N131").Value = -34150.9395
$ws.Range("H132").Value = 1174.7693
$ws.Range("I132").Value = 975
$ws.Range("J132").Value = 1346
$ws.Range("K132").Value = 8775
$ws.Range("L132").Value = 12114
$ws.Range("M132").Value = -6245
$ws.Range("N132").Value = -17174
$ws.Range("H135").Value = 915.55884
$ws.Range("I135").Value = 617.8421
$ws.Range("J135").Value = 1292.6666
$ws.Range("K135").Value = 5560.5789
$ws.Range("L135").Value = 11633.9994
$ws.Range("M135").Value = -3025.5789
$ws.Range("N135").Value = -16703.9994

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1059.95
$ws.Range("I22").Value = 966.5
$ws.Range("J22").Value = 1100
$ws.Range("K22").Value = 966.5
$ws.Range("L22").Value = 1100
$ws.Range("M22").Value = -671.5
$ws.Range("N22").Value = -1690
$ws.Range("H27").Value = 1059.95
$ws.Range("I27").Value = 966.5
$ws.Range("J27").Value = 1100
$ws.Range("K27").Value = 966.5
$ws.Range("L27").Value = 1100
$ws.Range("M27").Value = -859.5
$ws.Range("N27").Value = -1314
$ws.Range("H132").Value = 3439.8708
$ws.Range("I132").Value = 3645.1304
$ws.Range("J132").Value = 2849.75
$ws.Range("K132").Value = 10935.3912
$ws.Range("L132").Value = 8549.25
$ws.Range("M132").Value = -8405.3912
$ws.Range("N132").Value = -13609.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1331.7826
$ws.Range("I126").Value = 1349.8667
$ws.Range("J126").Value = 1297.875
$ws.Range("K126").Value = 4049.6001
$ws.Range("L126").Value = 3893.625
$ws.Range("M126").Value = -1579.6001
$ws.Range("N126").Value = -8833.625

